# Add a new "Slovakia" market sheet, cloned from the "Portugal" sheet
# (same layout/styles/merges/column widths), then update its market-specific
# values and selection state. Mirrors: "Test Data Added for Slovakia market"

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Clone Portugal -> new sheet is inserted right after it and becomes active.
[void]$portugal.Copy($null, $portugal)
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# The cloned sheet inherited Portugal's custom row heights (28.8) on rows
# 3-5; the Slovakia sheet uses default row heights, so reset them.
[void]$slovakia.Rows("3:5").AutoFit()

# Fill in the market-specific values. Order matters so new shared strings
# land in the same slots as the source data (NGC code before market name).
$slovakia.Range("B4").Value = "NGC-2930/T3178"
$slovakia.Range("B2").Value = "Slovakia Market"

# Match final UI state: Slovakia tab selected with A8 active, Portugal's
# selection left on the full sheet (no single active cell).
[void]$portugal.Cells.Select()
$slovakia.Activate()
[void]$slovakia.Range("A8").Select()
